# Update the localization-status report to reflect the latest handoff run.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-34-20 14:34:29"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-20 14:34:26"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-20 14:34:29"
